# Update auto scs (lamda_1, col B), time in ms (lamda_2, col C) and the
# auto-capacity Poisson dictionary (cols D/E) for rows 2-55, and drop the
# now-obsolete last row (56) so the table ends at row 55.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = New-Object 'object[,]' 54,4
$newValues[0,0] = 33.94444444444444; $newValues[0,1] = 1.95; $newValues[0,2] = 0; $newValues[0,3] = 0.148
$newValues[1,0] = 33.94444444444444; $newValues[1,1] = 1.95; $newValues[1,2] = 3; $newValues[1,3] = 0.006
$newValues[2,0] = 33.94444444444444; $newValues[2,1] = 1.95; $newValues[2,2] = 4; $newValues[2,3] = 0.013
$newValues[3,0] = 33.94444444444444; $newValues[3,1] = 1.95; $newValues[3,2] = 5; $newValues[3,3] = 0.019
$newValues[4,0] = 33.94444444444444; $newValues[4,1] = 1.95; $newValues[4,2] = 6; $newValues[4,3] = 0.039
$newValues[5,0] = 33.94444444444444; $newValues[5,1] = 1.95; $newValues[5,2] = 7; $newValues[5,3] = 0.05
$newValues[6,0] = 33.94444444444444; $newValues[6,1] = 1.95; $newValues[6,2] = 8; $newValues[6,3] = 0.044
$newValues[7,0] = 33.94444444444444; $newValues[7,1] = 1.95; $newValues[7,2] = 9; $newValues[7,3] = 0.038
$newValues[8,0] = 33.94444444444444; $newValues[8,1] = 1.95; $newValues[8,2] = 10; $newValues[8,3] = 0.045
$newValues[9,0] = 33.94444444444444; $newValues[9,1] = 1.95; $newValues[9,2] = 11; $newValues[9,3] = 0.038
$newValues[10,0] = 33.94444444444444; $newValues[10,1] = 1.95; $newValues[10,2] = 12; $newValues[10,3] = 0.022
$newValues[11,0] = 33.94444444444444; $newValues[11,1] = 1.95; $newValues[11,2] = 13; $newValues[11,3] = 0.027
$newValues[12,0] = 33.94444444444444; $newValues[12,1] = 1.95; $newValues[12,2] = 14; $newValues[12,3] = 0.024
$newValues[13,0] = 33.94444444444444; $newValues[13,1] = 1.95; $newValues[13,2] = 15; $newValues[13,3] = 0.041
$newValues[14,0] = 33.94444444444444; $newValues[14,1] = 1.95; $newValues[14,2] = 16; $newValues[14,3] = 0.037
$newValues[15,0] = 33.94444444444444; $newValues[15,1] = 1.95; $newValues[15,2] = 17; $newValues[15,3] = 0.038
$newValues[16,0] = 33.94444444444444; $newValues[16,1] = 1.95; $newValues[16,2] = 18; $newValues[16,3] = 0.03
$newValues[17,0] = 33.94444444444444; $newValues[17,1] = 1.95; $newValues[17,2] = 19; $newValues[17,3] = 0.031
$newValues[18,0] = 33.94444444444444; $newValues[18,1] = 1.95; $newValues[18,2] = 20; $newValues[18,3] = 0.03
$newValues[19,0] = 33.94444444444444; $newValues[19,1] = 1.95; $newValues[19,2] = 21; $newValues[19,3] = 0.022
$newValues[20,0] = 33.94444444444444; $newValues[20,1] = 1.95; $newValues[20,2] = 22; $newValues[20,3] = 0.017
$newValues[21,0] = 33.94444444444444; $newValues[21,1] = 1.95; $newValues[21,2] = 23; $newValues[21,3] = 0.016
$newValues[22,0] = 33.94444444444444; $newValues[22,1] = 1.95; $newValues[22,2] = 24; $newValues[22,3] = 0.013
$newValues[23,0] = 33.94444444444444; $newValues[23,1] = 1.95; $newValues[23,2] = 25; $newValues[23,3] = 0.022
$newValues[24,0] = 33.94444444444444; $newValues[24,1] = 1.95; $newValues[24,2] = 26; $newValues[24,3] = 0.02
$newValues[25,0] = 33.94444444444444; $newValues[25,1] = 1.95; $newValues[25,2] = 27; $newValues[25,3] = 0.019
$newValues[26,0] = 33.94444444444444; $newValues[26,1] = 1.95; $newValues[26,2] = 28; $newValues[26,3] = 0.021
$newValues[27,0] = 33.94444444444444; $newValues[27,1] = 1.95; $newValues[27,2] = 29; $newValues[27,3] = 0.014
$newValues[28,0] = 33.94444444444444; $newValues[28,1] = 1.95; $newValues[28,2] = 30; $newValues[28,3] = 0.007
$newValues[29,0] = 33.94444444444444; $newValues[29,1] = 1.95; $newValues[29,2] = 31; $newValues[29,3] = 0.014
$newValues[30,0] = 33.94444444444444; $newValues[30,1] = 1.95; $newValues[30,2] = 32; $newValues[30,3] = 0.009000000000000001
$newValues[31,0] = 33.94444444444444; $newValues[31,1] = 1.95; $newValues[31,2] = 33; $newValues[31,3] = 0.006
$newValues[32,0] = 33.94444444444444; $newValues[32,1] = 1.95; $newValues[32,2] = 34; $newValues[32,3] = 0.007
$newValues[33,0] = 33.94444444444444; $newValues[33,1] = 1.95; $newValues[33,2] = 35; $newValues[33,3] = 0.011
$newValues[34,0] = 33.94444444444444; $newValues[34,1] = 1.95; $newValues[34,2] = 36; $newValues[34,3] = 0.006
$newValues[35,0] = 33.94444444444444; $newValues[35,1] = 1.95; $newValues[35,2] = 37; $newValues[35,3] = 0.007
$newValues[36,0] = 33.94444444444444; $newValues[36,1] = 1.95; $newValues[36,2] = 38; $newValues[36,3] = 0.008
$newValues[37,0] = 33.94444444444444; $newValues[37,1] = 1.95; $newValues[37,2] = 39; $newValues[37,3] = 0.007
$newValues[38,0] = 33.94444444444444; $newValues[38,1] = 1.95; $newValues[38,2] = 40; $newValues[38,3] = 0.005
$newValues[39,0] = 33.94444444444444; $newValues[39,1] = 1.95; $newValues[39,2] = 41; $newValues[39,3] = 0.004
$newValues[40,0] = 33.94444444444444; $newValues[40,1] = 1.95; $newValues[40,2] = 42; $newValues[40,3] = 0.005
$newValues[41,0] = 33.94444444444444; $newValues[41,1] = 1.95; $newValues[41,2] = 43; $newValues[41,3] = 0.003
$newValues[42,0] = 33.94444444444444; $newValues[42,1] = 1.95; $newValues[42,2] = 44; $newValues[42,3] = 0.002
$newValues[43,0] = 33.94444444444444; $newValues[43,1] = 1.95; $newValues[43,2] = 45; $newValues[43,3] = 0.001
$newValues[44,0] = 33.94444444444444; $newValues[44,1] = 1.95; $newValues[44,2] = 46; $newValues[44,3] = 0.001
$newValues[45,0] = 33.94444444444444; $newValues[45,1] = 1.95; $newValues[45,2] = 47; $newValues[45,3] = 0.002
$newValues[46,0] = 33.94444444444444; $newValues[46,1] = 1.95; $newValues[46,2] = 48; $newValues[46,3] = 0.002
$newValues[47,0] = 33.94444444444444; $newValues[47,1] = 1.95; $newValues[47,2] = 49; $newValues[47,3] = 0.001
$newValues[48,0] = 33.94444444444444; $newValues[48,1] = 1.95; $newValues[48,2] = 52; $newValues[48,3] = 0.001
$newValues[49,0] = 33.94444444444444; $newValues[49,1] = 1.95; $newValues[49,2] = 54; $newValues[49,3] = 0.002
$newValues[50,0] = 33.94444444444444; $newValues[50,1] = 1.95; $newValues[50,2] = 57; $newValues[50,3] = 0.001
$newValues[51,0] = 33.94444444444444; $newValues[51,1] = 1.95; $newValues[51,2] = 61; $newValues[51,3] = 0.001
$newValues[52,0] = 33.94444444444444; $newValues[52,1] = 1.95; $newValues[52,2] = 63; $newValues[52,3] = 0.001
$newValues[53,0] = 33.94444444444444; $newValues[53,1] = 1.95; $newValues[53,2] = 67; $newValues[53,3] = 0.001

$ws.Range("B2:E55").Value = $newValues

# Row 56 (old last entry) no longer exists in the updated table.
$ws.Rows.Item(56).Delete()
